# SwaadSutra_Daily_2026-01-13.xlsx update
# New incoming order (#5, Sagar Borse - Til Poli x1) is inserted at the top of the
# "Daily Orders" log, pushing the previously logged orders down by one row. The
# "Summary" and "Items Breakdown" sheets are refreshed accordingly.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: "Daily Orders" - insert the new order as row 2, shifting the rest down
# ---------------------------------------------------------------------------
$orders = $wb.Worksheets.Item("Daily Orders")

$orders.Rows.Item(2).Insert()

$orders.Cells.Item(2, 1).Value = 5
$orders.Cells.Item(2, 2).NumberFormat = "@"
$orders.Cells.Item(2, 2).Value = "2026-01-13 16:40"
$orders.Cells.Item(2, 3).Value = "Sagar Borse"
$orders.Cells.Item(2, 4).Value = "A-1608"
$orders.Cells.Item(2, 5).NumberFormat = "@"
$orders.Cells.Item(2, 5).Value = "7588930329"
$orders.Cells.Item(2, 6).Value = "Til Poli x1"
$orders.Cells.Item(2, 7).Value = 30
$orders.Cells.Item(2, 8).Value = "NEW"
$orders.Cells.Item(2, 9).Value = "PENDING"
$orders.Cells.Item(2, 10).NumberFormat = "@"
$orders.Cells.Item(2, 10).Value = "2026-01-14"
$orders.Cells.Item(2, 11).NumberFormat = "@"
$orders.Cells.Item(2, 11).Value = "10:00"

# ---------------------------------------------------------------------------
# Sheet 2: "Summary" - refresh the aggregate counters
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")

$summary.Cells.Item(2, 1).Value = 5    # Total Orders
$summary.Cells.Item(2, 2).Value = 4    # New
$summary.Cells.Item(2, 3).Value = 1    # Cooking
$summary.Cells.Item(2, 4).Value = 0    # Ready
$summary.Cells.Item(2, 5).Value = 0    # Delivered
$summary.Cells.Item(2, 6).Value = 0    # Cancelled
$summary.Cells.Item(2, 7).Value = 165  # Total Revenue
$summary.Cells.Item(2, 8).Value = 0    # Paid Amount

# ---------------------------------------------------------------------------
# Sheet 3: "Items Breakdown" - Til Poli now leads with the extra order, so the
# rows are re-ranked by quantity
# ---------------------------------------------------------------------------
$items = $wb.Worksheets.Item("Items Breakdown")

$items.Cells.Item(2, 1).Value = "Til Poli"
$items.Cells.Item(2, 2).Value = 2
$items.Cells.Item(2, 3).Value = 60

$items.Cells.Item(3, 1).Value = "Onion Pakoda (Kanda Bhaje)"
$items.Cells.Item(3, 2).Value = 1
$items.Cells.Item(3, 3).Value = 60

$items.Cells.Item(4, 1).Value = "Pohe"
$items.Cells.Item(4, 2).Value = 1
$items.Cells.Item(4, 3).Value = 30

$items.Cells.Item(5, 1).Value = "Wheat Chapati"
$items.Cells.Item(5, 2).Value = 1
$items.Cells.Item(5, 3).Value = 15
